{"js": "// Update the date line and each \"a\u00f7b=\" problem text to the new values\n// described by the commit diff. Every \"old\" string is unique in the\n// document, so a simple matchCase search-and-replace per pair is safe.\nconst replacements = [\n  [\"2024-08-06 Tuesday\", \"2024-08-07 Wednesday\"],\n  [\"32\u00f79=\", \"36\u00f74=\"],\n  [\"62\u00f77=\", \"38\u00f72=\"],\n  [\"50\u00f79=\", \"12\u00f72=\"],\n  [\"12\u00f77=\", \"91\u00f77=\"],\n  [\"17\u00f77=\", \"90\u00f79=\"],\n  [\"44\u00f76=\", \"64\u00f74=\"],\n  [\"19\u00f77=\", \"62\u00f74=\"],\n  [\"21\u00f79=\", \"80\u00f73=\"],\n  [\"97\u00f73=\", \"77\u00f78=\"],\n  [\"30\u00f74=\", \"34\u00f77=\"],\n  [\"77\u00f72=\", \"25\u00f75=\"],\n  [\"16\u00f76=\", \"58\u00f73=\"],\n  [\"88\u00f75=\", \"68\u00f73=\"],\n  [\"84\u00f74=\", \"79\u00f72=\"],\n  [\"26\u00f72=\", \"89\u00f72=\"],\n  [\"24\u00f72=\", \"10\u00f78=\"],\n  [\"45\u00f77=\", \"24\u00f74=\"],\n  [\"46\u00f77=\", \"24\u00f79=\"],\n  [\"64\u00f76=\", \"28\u00f79=\"],\n  [\"79\u00f76=\", \"48\u00f75=\"],\n  [\"96\u00f72=\", \"14\u00f75=\"],\n  [\"50\u00f78=\", \"25\u00f75=\"],\n  [\"71\u00f74=\", \"16\u00f78=\"],\n  [\"54\u00f77=\", \"28\u00f74=\"],\n  [\"52\u00f73=\", \"73\u00f78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-06 Tuesday\", \"2024-08-07 Wednesday\"),\n    @(\"32\u00f79=\", \"36\u00f74=\"),\n    @(\"62\u00f77=\", \"38\u00f72=\"),\n    @(\"50\u00f79=\", \"12\u00f72=\"),\n    @(\"12\u00f77=\", \"91\u00f77=\"),\n    @(\"17\u00f77=\", \"90\u00f79=\"),\n    @(\"44\u00f76=\", \"64\u00f74=\"),\n    @(\"19\u00f77=\", \"62\u00f74=\"),\n    @(\"21\u00f79=\", \"80\u00f73=\"),\n    @(\"97\u00f73=\", \"77\u00f78=\"),\n    @(\"30\u00f74=\", \"34\u00f77=\"),\n    @(\"77\u00f72=\", \"25\u00f75=\"),\n    @(\"16\u00f76=\", \"58\u00f73=\"),\n    @(\"88\u00f75=\", \"68\u00f73=\"),\n    @(\"84\u00f74=\", \"79\u00f72=\"),\n    @(\"26\u00f72=\", \"89\u00f72=\"),\n    @(\"24\u00f72=\", \"10\u00f78=\"),\n    @(\"45\u00f77=\", \"24\u00f74=\"),\n    @(\"46\u00f77=\", \"24\u00f79=\"),\n    @(\"64\u00f76=\", \"28\u00f79=\"),\n    @(\"79\u00f76=\", \"48\u00f75=\"),\n    @(\"96\u00f72=\", \"14\u00f75=\"),\n    @(\"50\u00f78=\", \"25\u00f75=\"),\n    @(\"71\u00f74=\", \"16\u00f78=\"),\n    @(\"54\u00f77=\", \"28\u00f74=\"),\n    @(\"52\u00f73=\", \"73\u00f78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
